$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B1:B2 hold the seed values for the Fibonacci-style sequence
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 2

# B3:B10 share one formula (B1+B2) filled down, matching the shared formula
# group (t="shared") produced by Excel when a formula is entered into a
# multi-cell range in one go.
$ws.Range("B3:B10").Formula = "=B1+B2"

# Reproduce the saved selection state (active cell B3, selection B3:B10)
$ws.Range("B3:B10").Select() | Out-Null
